$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the speculative Who/Time entries that had not actually been submitted yet
$ws.Range("E2:F3").ClearContents()
$ws.Range("E10:F15").ClearContents()

# Mark "check if move is valid" (row 7) as finished, with who/time actually submitted
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "R"
$ws.Range("F7").Value = 180

# Update the active selection to reflect where the user was working
[void]$ws.Range("E10").Select()
